$wb = $excel.ActiveWorkbook

# Rename the "Contact" sheet to "Project - Contact" (added new nested tab name)
$ws = $wb.Worksheets.Item("Contact")
$ws.Name = "Project - Contact"
